# "moved second slide box over left"
#
# Every top-level shape on slide 2 (except the header/banner group
# "Group 51") is shifted left by 66504 EMU, keeping Top unchanged.
# This mirrors the author dragging/nudging the whole "box" of content
# on the second slide to the left.
#
# PowerPoint's COM object model expresses Shape.Left/Top in points
# (1 pt = 12700 EMU), while the underlying OOXML stores EMU. To land on
# the exact target EMU values from the diff despite the runtime
# truncating (not rounding) points->EMU internally, we add a tiny
# safety epsilon (well under half an EMU in point-terms) before
# assigning, which reliably lands on the desired integer EMU without
# overshooting to the next one.

$EMU_PER_POINT = 12700
$EPS = 0.00003

# Target absolute Left, in EMU, for every affected top-level shape on
# slide 2 (taken from the canonical OOXML diff).
$targets = @{
    "Rectangle 38" = 288994
    "TextBox 39"   = 3036746
    "TextBox 41"   = 477478
    "TextBox 42"   = 477478
    "TextBox 43"   = 6029496
    "TextBox 44"   = 477479
    "Rectangle 45" = 6138352
    "TextBox 47"   = 477478
    "Group 49"     = 534853
    "TextBox 50"   = 477478
    "Rectangle 52" = 534852
    "Oval 54"      = 635596
    "Oval 56"      = 635596
    "Oval 57"      = 644927
    "TextBox 58"   = 839355
    "TextBox 59"   = 839355
    "TextBox 60"   = 839354
    "Rectangle 61" = 3255194
    "Rectangle 62" = 713754
    "Rectangle 63" = 713754
    "TextBox 65"   = 1204845
    "Graphic 68"   = 721815
    "TextBox 69"   = 1204845
    "TextBox 70"   = 6029496
    "Rectangle 71" = 6138352
    "TextBox 73"   = 9586256
    "Group 3"      = 11742001
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    $name = $shp.Name
    if ($targets.ContainsKey($name)) {
        $targetEmu = $targets[$name]
        $shp.Left = ($targetEmu / $EMU_PER_POINT) + $EPS
    }
}
